$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row: two extra columns (ownTeam, oppTeam) inserted before "batsman",
# pushing batsman/totalRuns/totalBalls/total4s/total6s/sr from D..I to F..K.
$headers = @("venue","date","result","ownTeam","oppTeam","batsman","totalRuns","totalBalls","total4s","total6s","sr")

# Full replacement data set (15 rows, reordered + 5 new match rows), columns A..K.
$data = @()
$data += ,@(" Abu Dhabi"," October 28 2020","Mumbai won by 5 wickets (with 5 balls remaining)","Royal Challengers Bangalore","Mumbai Indians","Devdutt Padikkal ","74","45","12","1","164.44")
$data += ,@(" Dubai (DSC)"," October 17 2020","RCB won by 7 wickets (with 2 balls remaining)","Royal Challengers Bangalore","Rajasthan Royals","Devdutt Padikkal ","35","37","2","0","94.59")
$data += ,@(" Abu Dhabi"," October 03 2020","RCB won by 8 wickets (with 5 balls remaining)","Royal Challengers Bangalore","Rajasthan Royals","Devdutt Padikkal ","63","45","6","1","140.00")
$data += ,@(" Dubai (DSC)"," September 24 2020","Kings XI won by 97 runs","Royal Challengers Bangalore","Kings XI Punjab","Devdutt Padikkal ","1","2","0","0","50.00")
$data += ,@(" Sharjah"," October 15 2020","Kings XI won by 8 wickets","Royal Challengers Bangalore","Kings XI Punjab","Devdutt Padikkal ","18","12","1","1","150.00")
$data += ,@(" Dubai (DSC)"," September 21 2020","RCB won by 10 runs","Royal Challengers Bangalore","Sunrisers Hyderabad","Devdutt Padikkal ","56","42","8","0","133.33")
$data += ,@(" Abu Dhabi"," November 06 2020","Sunrisers won by 6 wickets (with 2 balls remaining)","Royal Challengers Bangalore","Sunrisers Hyderabad","Devdutt Padikkal ","1","6","0","0","16.66")
$data += ,@(" Dubai (DSC)"," October 05 2020","Capitals won by 59 runs","Royal Challengers Bangalore","Delhi Capitals","Devdutt Padikkal ","4","6","0","0","66.66")
$data += ,@(" Sharjah"," October 31 2020","Sunrisers won by 5 wickets (with 35 balls remaining)","Royal Challengers Bangalore","Sunrisers Hyderabad","Devdutt Padikkal ","5","8","1","0","62.50")
$data += ,@(" Abu Dhabi"," November 02 2020","Capitals won by 6 wickets (with 6 balls remaining)","Royal Challengers Bangalore","Delhi Capitals","Devdutt Padikkal ","50","41","5","0","121.95")
$data += ,@(" Abu Dhabi"," October 21 2020","RCB won by 8 wickets (with 39 balls remaining)","Royal Challengers Bangalore","Kolkata Knight Riders","Devdutt Padikkal ","25","17","3","0","147.05")
$data += ,@(" Dubai (DSC)"," September 28 2020","Match tied (RCB won the one-over eliminator)","Royal Challengers Bangalore","Mumbai Indians","Devdutt Padikkal ","54","40","5","2","135.00")
$data += ,@(" Sharjah"," October 12 2020","RCB won by 82 runs","Royal Challengers Bangalore","Kolkata Knight Riders","Devdutt Padikkal ","32","23","4","1","139.13")
$data += ,@(" Dubai (DSC)"," October 25 2020","Super Kings won by 8 wickets (with 8 balls remaining)","Royal Challengers Bangalore","Chennai Super Kings","Devdutt Padikkal ","22","21","2","1","104.76")
$data += ,@(" Dubai (DSC)"," October 10 2020","RCB won by 37 runs","Royal Challengers Bangalore","Chennai Super Kings","Devdutt Padikkal ","33","34","2","1","97.05")

# Pre-format the full target range as Text so every value (including digit-only
# strings like "74" or "0") is stored as text, matching the source data which
# keeps numbers-as-text (see the sheet's numberStoredAsText ignoredError).
$fullRange = $ws.Range("A1:K16")
$fullRange.NumberFormat = "@"

# Header row
for ($c = 1; $c -le $headers.Count; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Data rows (start at worksheet row 2)
for ($r = 0; $r -lt $data.Count; $r++) {
    $rowVals = $data[$r]
    $excelRow = $r + 2
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws.Cells.Item($excelRow, $c).Value = $rowVals[$c - 1]
    }
}
